# Trade #18 closed at 2026-02-17 08:19:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")

# --- Summary sheet updates ---
$wsSummary.Range("B3").Value = 1200
$wsSummary.Range("B4").Value = 0
$wsSummary.Range("B5").Value = 0
$wsSummary.Range("B6").Value = 18
$wsSummary.Range("B8").Value = 9
$wsSummary.Range("B9").Value = 27.78

# --- Strategy Status sheet updates (MarketMaking row, row 4) ---
$wsStrategy.Range("C4").Value = 100
$wsStrategy.Range("D4").Value = 18
$wsStrategy.Range("E4").Value = 0
$wsStrategy.Range("F4").Value = -0
$wsStrategy.Range("G4").Value = 27.78

# --- Add new trade row (#18) to "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow18($ws) {
    $ws.Cells.Item(19, 1).Value = 18
    # Force the date-looking string to stay plain text instead of being
    # auto-converted into a date serial number by the smart-entry parser.
    $ws.Cells.Item(19, 2).NumberFormat = "@"
    $ws.Cells.Item(19, 2).Value = "2026-02-17"
    $ws.Cells.Item(19, 2).Style = "Normal"
    $ws.Cells.Item(19, 3).Value = "08:19:29"
    $ws.Cells.Item(19, 4).Value = "MarketMaking"
    $ws.Cells.Item(19, 5).Value = "UP"
    $ws.Cells.Item(19, 6).Value = 0.03
    $ws.Cells.Item(19, 7).Value = 0.02
    $ws.Cells.Item(19, 8).Value = "CLOSED"
    $ws.Cells.Item(19, 9).Value = -33.3333
    $ws.Cells.Item(19, 10).Value = -0.01
    $ws.Cells.Item(19, 11).Value = 100
    $ws.Cells.Item(19, 12).Value = 0
    $ws.Cells.Item(19, 13).Value = 0
    $ws.Cells.Item(19, 14).Value = 0.6
    $ws.Cells.Item(19, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(19, 16).Value = "early_exit"
    $ws.Cells.Item(19, 17).Value = 0.13
}

Add-TradeRow18 $wsAllTrades
Add-TradeRow18 $wsMarketMaking
